# Implemented BruteForce Traceback solution
# Fill in the Sudoku puzzle cells that the brute-force/backtracking solver produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New numbers placed by the solver
$ws.Range("G1").Value = 1
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 3

$ws.Range("A2").Value = 8
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 4
$ws.Range("I2").Value = 2

$ws.Range("A3").Value = 2

$ws.Range("B4").Value = 7
$ws.Range("E4").Value = 5

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 8
$ws.Range("G5").ClearContents()
$ws.Range("I5").Value = 9

$ws.Range("C6").Value = 8
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 9

$ws.Range("C7").Value = 5
$ws.Range("E7").Value = 9
$ws.Range("H7").Value = 3

$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 3
$ws.Range("G8").Value = 5
$ws.Range("I8").Value = 7

$ws.Range("F9").Value = 4
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = 8
$ws.Range("I9").ClearContents()

# Window / view state changes
$ws.Range("M13").Select()
$excel.ActiveWindow.Zoom = 190
